# Generate Report for Handback
# Update the handoff/handback timestamp cells on each worksheet.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the first data row
$overview.Range("G2").Value = "2017-01-03 08:10:25"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (L2)
$zhcn.Range("H2").Value = "2017-01-03 08:10:14"
$zhcn.Range("L2").Value = "2017-01-03 08:10:50"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (L2)
$dede.Range("H2").Value = "2017-01-03 08:10:25"
$dede.Range("L2").Value = "2017-01-03 08:11:03"
